$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.293025333333333
$ws.Range("H2").Value = 9.879076
$ws.Range("I2").Value = 0.2347004583311417
$ws.Range("J2").Value = 0.2347004583311417
$ws.Range("M2").Value = 8.540560666666666
$ws.Range("N2").Value = 25.621682
$ws.Range("O2").Value = 0.4159358086620884
$ws.Range("P2").Value = 0.4159358086620884
$ws.Range("Q2").Value = 28.12428263620355
$ws.Range("R2").Value = 253.118543725832
$ws.Range("S2").Value = 0.09762032492932622
$ws.Range("T2").Value = 0.09762032492932622

$ws.Range("G3").Value = 3.293025333333333
$ws.Range("H3").Value = 9.879076
$ws.Range("I3").Value = 0.2347004583311417
$ws.Range("J3").Value = 0.2347004583311417
$ws.Range("O3").Value = 0.563694901924408
$ws.Range("P3").Value = 0.563694901924408
$ws.Range("Q3").Value = 38.11529186030889
$ws.Range("R3").Value = 343.03762674278
$ws.Range("S3").Value = 0.1322994518405865
$ws.Range("T3").Value = 0.1322994518405865

$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("G4").Value = 3.293025333333333
$ws.Range("H4").Value = 9.879076
$ws.Range("I4").Value = 0.2347004583311417
$ws.Range("J4").Value = 0.2347004583311417
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.026642
$ws.Range("N4").Value = 0.079926
$ws.Range("O4").Value = 0.001297498167494471
$ws.Range("P4").Value = 0.001297498167494471
$ws.Range("Q4").Value = 0.08773278093066665
$ws.Range("R4").Value = 0.7895950283759999
$ws.Range("S4").Value = 0.0003045234145947689
$ws.Range("T4").Value = 0.0003045234145947689

$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 3.293025333333333
$ws.Range("H5").Value = 9.879076
$ws.Range("I5").Value = 0.2347004583311417
$ws.Range("J5").Value = 0.2347004583311417
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.391608
$ws.Range("N5").Value = 1.174824
$ws.Range("O5").Value = 0.01907179124600912
$ws.Range("P5").Value = 0.01907179124600912
$ws.Range("Q5").Value = 1.289575064736
$ws.Range("R5").Value = 11.606175582624
$ws.Range("S5").Value = 0.004476158146634196
$ws.Range("T5").Value = 0.004476158146634196

$ws.Range("I6").Value = 0.422318927221656
$ws.Range("J6").Value = 0.422318927221656
$ws.Range("M6").Value = 8.540560666666666
$ws.Range("N6").Value = 25.621682
$ws.Range("O6").Value = 0.4159358086620884
$ws.Range("P6").Value = 0.4159358086620884
$ws.Range("Q6").Value = 50.60670505824978
$ws.Range("R6").Value = 455.460345524248
$ws.Range("S6").Value = 0.1756575645072452
$ws.Range("T6").Value = 0.1756575645072452

$ws.Range("I7").Value = 0.422318927221656
$ws.Range("J7").Value = 0.422318927221656
$ws.Range("O7").Value = 0.563694901924408
$ws.Range("P7").Value = 0.563694901924408
$ws.Range("S7").Value = 0.2380590262610326
$ws.Range("T7").Value = 0.2380590262610326

$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("I8").Value = 0.422318927221656
$ws.Range("J8").Value = 0.422318927221656
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.026642
$ws.Range("N8").Value = 0.079926
$ws.Range("O8").Value = 0.001297498167494471
$ws.Range("P8").Value = 0.001297498167494471
$ws.Range("Q8").Value = 0.1578659632293333
$ws.Range("R8").Value = 1.420793669064
$ws.Range("S8").Value = 0.0005479580341683296
$ws.Range("T8").Value = 0.0005479580341683296

$ws.Range("D9").Value = "MuSCs"
$ws.Range("I9").Value = 0.422318927221656
$ws.Range("J9").Value = 0.422318927221656
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.391608
$ws.Range("N9").Value = 1.174824
$ws.Range("O9").Value = 0.01907179124600912
$ws.Range("P9").Value = 0.01907179124600912
$ws.Range("Q9").Value = 2.320455451104
$ws.Range("R9").Value = 20.884099059936
$ws.Range("S9").Value = 0.00805437841920994
$ws.Range("T9").Value = 0.00805437841920994

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.186484333333333
$ws.Range("H10").Value = 6.559453
$ws.Range("I10").Value = 0.1558350827042511
$ws.Range("J10").Value = 0.1558350827042511
$ws.Range("M10").Value = 8.540560666666666
$ws.Range("N10").Value = 25.621682
$ws.Range("O10").Value = 0.4159358086620884
$ws.Range("P10").Value = 0.4159358086620884
$ws.Range("Q10").Value = 18.67380209554955
$ws.Range("R10").Value = 168.064218859946
$ws.Range("S10").Value = 0.06481739114251613
$ws.Range("T10").Value = 0.06481739114251613

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 2.186484333333333
$ws.Range("H11").Value = 6.559453
$ws.Range("I11").Value = 0.1558350827042511
$ws.Range("J11").Value = 0.1558350827042511
$ws.Range("O11").Value = 0.563694901924408
$ws.Range("P11").Value = 0.563694901924408
$ws.Range("Q11").Value = 25.30757588452389
$ws.Range("R11").Value = 227.768182960715
$ws.Range("S11").Value = 0.08784344166135485
$ws.Range("T11").Value = 0.08784344166135485

$ws.Range("D12").Value = "Inflammatory-Mac"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 2.186484333333333
$ws.Range("H12").Value = 6.559453
$ws.Range("I12").Value = 0.1558350827042511
$ws.Range("J12").Value = 0.1558350827042511
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.026642
$ws.Range("N12").Value = 0.079926
$ws.Range("O12").Value = 0.001297498167494471
$ws.Range("P12").Value = 0.001297498167494471
$ws.Range("Q12").Value = 0.05825231560866666
$ws.Range("R12").Value = 0.524270840478
$ws.Range("S12").Value = 0.0002021957342401152
$ws.Range("T12").Value = 0.0002021957342401152

$ws.Range("D13").Value = "MuSCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 2.186484333333333
$ws.Range("H13").Value = 6.559453
$ws.Range("I13").Value = 0.1558350827042511
$ws.Range("J13").Value = 0.1558350827042511
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.391608
$ws.Range("N13").Value = 1.174824
$ws.Range("O13").Value = 0.01907179124600912
$ws.Range("P13").Value = 0.01907179124600912
$ws.Range("Q13").Value = 0.8562447568079999
$ws.Range("R13").Value = 7.706202811272
$ws.Range("S13").Value = 0.002972054166140044
$ws.Range("T13").Value = 0.002972054166140044

$ws.Range("G14").Value = 1.169741
$ws.Range("H14").Value = 3.509223
$ws.Range("I14").Value = 0.08336976519729013
$ws.Range("J14").Value = 0.08336976519729013
$ws.Range("M14").Value = 8.540560666666666
$ws.Range("N14").Value = 25.621682
$ws.Range("O14").Value = 0.4159358086620884
$ws.Range("P14").Value = 0.4159358086620884
$ws.Range("Q14").Value = 9.990243974787331
$ws.Range("R14").Value = 89.912195773086
$ws.Range("S14").Value = 0.0346764707053033
$ws.Range("T14").Value = 0.0346764707053033

$ws.Range("G15").Value = 1.169741
$ws.Range("H15").Value = 3.509223
$ws.Range("I15").Value = 0.08336976519729013
$ws.Range("J15").Value = 0.08336976519729013
$ws.Range("O15").Value = 0.563694901924408
$ws.Range("P15").Value = 0.563694901924408
$ws.Range("Q15").Value = 13.53922764111833
$ws.Range("R15").Value = 121.853048770065
$ws.Range("S15").Value = 0.04699511161634738
$ws.Range("T15").Value = 0.04699511161634738

$ws.Range("D16").Value = "Inflammatory-Mac"
$ws.Range("G16").Value = 1.169741
$ws.Range("H16").Value = 3.509223
$ws.Range("I16").Value = 0.08336976519729013
$ws.Range("J16").Value = 0.08336976519729013
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.026642
$ws.Range("N16").Value = 0.079926
$ws.Range("O16").Value = 0.001297498167494471
$ws.Range("P16").Value = 0.001297498167494471
$ws.Range("Q16").Value = 0.031164239722
$ws.Range("R16").Value = 0.280478157498
$ws.Range("S16").Value = 0.0001081721175679283
$ws.Range("T16").Value = 0.0001081721175679283

$ws.Range("D17").Value = "MuSCs"
$ws.Range("G17").Value = 1.169741
$ws.Range("H17").Value = 3.509223
$ws.Range("I17").Value = 0.08336976519729013
$ws.Range("J17").Value = 0.08336976519729013
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.391608
$ws.Range("N17").Value = 1.174824
$ws.Range("O17").Value = 0.01907179124600912
$ws.Range("P17").Value = 0.01907179124600912
$ws.Range("Q17").Value = 0.458079933528
$ws.Range("R17").Value = 4.122719401752001
$ws.Range("S17").Value = 0.001590010758071513
$ws.Range("T17").Value = 0.001590010758071513

$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 1.456052666666666
$ws.Range("H18").Value = 4.368157999999999
$ws.Range("I18").Value = 0.1037757665456611
$ws.Range("J18").Value = 0.1037757665456611
$ws.Range("M18").Value = 8.540560666666666
$ws.Range("N18").Value = 25.621682
$ws.Range("O18").Value = 0.4159358086620884
$ws.Range("P18").Value = 0.4159358086620884
$ws.Range("Q18").Value = 12.43550613352844
$ws.Range("R18").Value = 111.919555201756
$ws.Range("S18").Value = 0.04316405737769764
$ws.Range("T18").Value = 0.04316405737769764

$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 1.456052666666666
$ws.Range("H19").Value = 4.368157999999999
$ws.Range("I19").Value = 0.1037757665456611
$ws.Range("J19").Value = 0.1037757665456611
$ws.Range("O19").Value = 0.563694901924408
$ws.Range("P19").Value = 0.563694901924408
$ws.Range("Q19").Value = 16.85315681972111
$ws.Range("R19").Value = 151.67841137749
$ws.Range("S19").Value = 0.05849787054508668
$ws.Range("T19").Value = 0.05849787054508668

$ws.Range("D20").Value = "Inflammatory-Mac"
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 1.456052666666666
$ws.Range("H20").Value = 4.368157999999999
$ws.Range("I20").Value = 0.1037757665456611
$ws.Range("J20").Value = 0.1037757665456611
$ws.Range("K20").Value = 1
$ws.Range("L20").Value = 0.3333333333333333
$ws.Range("M20").Value = 0.026642
$ws.Range("N20").Value = 0.079926
$ws.Range("O20").Value = 0.001297498167494471
$ws.Range("P20").Value = 0.001297498167494471
$ws.Range("Q20").Value = 0.03879215514533332
$ws.Range("R20").Value = 0.3491293963079999
$ws.Range("S20").Value = 0.0001346488669233293
$ws.Range("T20").Value = 0.0001346488669233293

$ws.Range("D21").Value = "MuSCs"
$ws.Range("E21").Value = 3
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 1.456052666666666
$ws.Range("H21").Value = 4.368157999999999
$ws.Range("I21").Value = 0.1037757665456611
$ws.Range("J21").Value = 0.1037757665456611
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 0.391608
$ws.Range("N21").Value = 1.174824
$ws.Range("O21").Value = 0.01907179124600912
$ws.Range("P21").Value = 0.01907179124600912
$ws.Range("Q21").Value = 0.5702018726879999
$ws.Range("R21").Value = 5.131816854192
$ws.Range("S21").Value = 0.001979189755953425
$ws.Range("T21").Value = 0.001979189755953425
